$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# All target cells hold text (inline strings) in the original workbook, including
# numeric-looking price strings that must keep exact formatting (trailing zeros, etc.).
# Prefix with a literal apostrophe so Excel stores them as text, not as numbers.
$ws.Range("D2").Value = "'249.28"
$ws.Range("D3").Value = "'21.90"
$ws.Range("D4").Value = "'5.546"
$ws.Range("D5").Value = "'0.05649"
$ws.Range("D6").Value = "'6.459"
$ws.Range("D7").Value = "'0.8005"
$ws.Range("D8").Value = "'1.036"
$ws.Range("B9").Value = "'One"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01170"
$ws.Range("E9").Value = "'8OneONEBestin24h"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1439"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07297"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03124"
$ws.Range("E12").Value = "'11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02915"
$ws.Range("E13").Value = "'12BitrueCoinBTR"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09269"
$ws.Range("E14").Value = "'13BitMartTokenBMX"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001674"
$ws.Range("E15").Value = "'14BitForexTokenBF"
$ws.Range("B16").Value = "'MCDex"
$ws.Range("C16").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.212"
$ws.Range("E16").Value = "'15MCDexMCB"
$ws.Range("D17").Value = "'0.04741"
$ws.Range("D18").Value = "'0.006421"
$ws.Range("D19").Value = "'0.005075"
$ws.Range("D20").Value = "'0.001050"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D24").Value = "'2.089"
$ws.Range("D25").Value = "'0.3268"
$ws.Range("D27").Value = "'0.0003301"
$ws.Range("D40").Value = "'0.04153"
$ws.Range("D41").Value = "'0.006943"
$ws.Range("B42").Value = "'CEJI"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003501"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1044"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.009392"
$ws.Range("D47").Value = "'0.6802"
$ws.Range("D48").Value = "'0.01604"
$ws.Range("D49").Value = "'0.00002101"
